$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for rows 2-10
# from serial date 45204 (2023-10-05) to 45207 (2023-10-08)
foreach ($row in 2..10) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value = 45207
    }
}
